$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 07:45"

# Update Israel row (row 27)
$ws.Range("B27").Value = 268175
$ws.Range("C27").Value = 1400
$ws.Range("D27").Value = 201392
$ws.Range("E27").Value = 65064

# Update Uzbekistan row (row 59)
$ws.Range("B59").Value = 58859
$ws.Range("C59").Value = 247
$ws.Range("D59").Value = 55318
$ws.Range("E59").Value = 3058
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 483

# Update Kirguistan row (row 66)
$ws.Range("B66").Value = 47635
$ws.Range("C66").Value = 207
$ws.Range("D66").Value = 43521
$ws.Range("E66").Value = 3048

# Swap Montserrat / Islas Malvinas rows (215/216), along with their stats
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
